# Auto-generated edit script: updates Leve profit-calculation value cells
# across all 8 job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) to refresh
# cached marketboard prices / profit figures, per the scheduled-runner commit.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 269.16666
$ws.Range("I2").Value = 214.28572
$ws.Range("J2").Value = 346
$ws.Range("K2").Value = 214.28572
$ws.Range("L2").Value = 346
$ws.Range("M2").Value = -101.28572
$ws.Range("N2").Value = -572
$ws.Range("H132").Value = 2321.0715
$ws.Range("I132").Value = 2469.16
$ws.Range("J132").Value = 1087
$ws.Range("K132").Value = 7407.48
$ws.Range("L132").Value = 3261
$ws.Range("M132").Value = -4877.48
$ws.Range("N132").Value = -8321
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
$ws.Range("H138").Value = 15877215
$ws.Range("I138").Value = 2009.72
$ws.Range("J138").Value = 26321428
$ws.Range("K138").Value = 6029.16
$ws.Range("L138").Value = 78964284
$ws.Range("M138").Value = -889.1599999999999
$ws.Range("N138").Value = -78974564

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 50500.5
$ws.Range("I4").Value = 50500.5
$ws.Range("K4").Value = 50500.5
$ws.Range("M4").Value = -50384.5
$ws.Range("H32").Value = 5683.14
$ws.Range("I32").Value = 3866.525
$ws.Range("K32").Value = 3866.525
$ws.Range("M32").Value = -3579.525

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 14511.429
$ws.Range("I22").Value = 14511.429
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 14511.429
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -14338.429
$ws.Range("N22").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 350.36365
$ws.Range("I22").Value = 350.36365
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 350.36365
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -0.3636500000000069
$ws.Range("N22").ClearContents()
$ws.Range("H134").Value = 4934.933
$ws.Range("I134").Value = 5324.923
$ws.Range("J134").Value = 2400
$ws.Range("K134").Value = 15974.769
$ws.Range("L134").Value = 7200
$ws.Range("M134").Value = -13439.769
$ws.Range("N134").Value = -12270

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 338158.88
$ws.Range("I5").Value = 386.20834
$ws.Range("J5").Value = 627678.3
$ws.Range("K5").Value = 1158.62502
$ws.Range("L5").Value = 1883034.9
$ws.Range("M5").Value = -1046.62502
$ws.Range("N5").Value = -1883258.9
$ws.Range("H22").Value = 1150
$ws.Range("I22").Value = 1150
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 3450
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -3281
$ws.Range("N22").ClearContents()
$ws.Range("H27").Value = 1150
$ws.Range("I27").Value = 1150
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 3450
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = -3348
$ws.Range("N27").ClearContents()
$ws.Range("H34").Value = 1715
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 1715
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 5145
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value = -5313
$ws.Range("H39").Value = 7279.8184
$ws.Range("J39").Value = 2864.4285
$ws.Range("L39").Value = 8593.2855
$ws.Range("N39").Value = -9181.2855
$ws.Range("H55").Value = 4925
$ws.Range("J55").Value = 4925
$ws.Range("L55").Value = 14775
$ws.Range("N55").Value = -15129
$ws.Range("H64").Value = 1096.4814
$ws.Range("J64").Value = 1178.875
$ws.Range("L64").Value = 3536.625
$ws.Range("N64").Value = -4076.625
$ws.Range("H67").Value = 1096.4814
$ws.Range("J67").Value = 1178.875
$ws.Range("L67").Value = 3536.625
$ws.Range("N67").Value = -5408.625
$ws.Range("H68").Value = 191984.83
$ws.Range("I68").Value = 593535.25
$ws.Range("K68").Value = 1780605.75
$ws.Range("M68").Value = -1779794.75
$ws.Range("H71").Value = 191984.83
$ws.Range("I71").Value = 593535.25
$ws.Range("K71").Value = 5341817.25
$ws.Range("M71").Value = -5337761.25
$ws.Range("H122").Value = 580
$ws.Range("I122").Value = 392.48
$ws.Range("J122").Value = 1361.3334
$ws.Range("K122").Value = 3532.32
$ws.Range("L122").Value = 12252.0006
$ws.Range("M122").Value = -1082.32
$ws.Range("N122").Value = -17152.0006
$ws.Range("H135").Value = 338158.88
$ws.Range("I135").Value = 386.20834
$ws.Range("J135").Value = 627678.3
$ws.Range("K135").Value = 3475.87506
$ws.Range("L135").Value = 5649104.7
$ws.Range("M135").Value = -940.8750600000003
$ws.Range("N135").Value = -5654174.7

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H41").Value = 3000
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 3000
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 3000
$ws.Range("M41").ClearContents()
$ws.Range("N41").Value = -3710
$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("M43").ClearContents()
$ws.Range("H70").Value = 28212
$ws.Range("I70").Value = 28212
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 28212
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -27942
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 28212
$ws.Range("I73").Value = 28212
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 28212
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -27276
$ws.Range("N73").ClearContents()
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("I22").Value = 709.1818
$ws.Range("J22").Value = 1027.1111
$ws.Range("K22").Value = 709.1818
$ws.Range("L22").Value = 1027.1111
$ws.Range("M22").Value = -414.1818
$ws.Range("N22").Value = -1617.1111
$ws.Range("I27").Value = 709.1818
$ws.Range("J27").Value = 1027.1111
$ws.Range("K27").Value = 709.1818
$ws.Range("L27").Value = 1027.1111
$ws.Range("M27").Value = -602.1818
$ws.Range("N27").Value = -1241.1111
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()
$ws.Range("H132").Value = 4106.05
$ws.Range("I132").Value = 4532.6772
$ws.Range("J132").Value = 3650
$ws.Range("K132").Value = 13598.0316
$ws.Range("L132").Value = 10950
$ws.Range("M132").Value = -11068.0316
$ws.Range("N132").Value = -16010

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 7941157.5
$ws.Range("I136").Value = 16667425
$ws.Range("J136").Value = 8187.5
$ws.Range("K136").Value = 50002275
$ws.Range("L136").Value = 24562.5
$ws.Range("M136").Value = -49999725
$ws.Range("N136").Value = -29662.5
